$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value2 = "https://www.youtube.com/watch?v=l-kxBjNML5c"
$ws.Range("A4").Value2 = "Cheerleading Fail "
$ws.Range("D4").Value2 = $ws.Range("D2").Value2
$ws.Range("E4").Formula = "=MID(F4, 33, 20)"
$ws.Range("C4").Formula = "=D4&MID(F4, 33, 20)&`$G`$2"
$ws.Range("B4").Value2 = $ws.Range("C4").Value2

$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F4").Value2 = "https://www.youtube.com/watch?v=l-kxBjNML5c"

$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.youtube.com/watch?v=l-kxBjNML5c")

$ws.Range("B4").Select()
